$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (Fri Dec 15 10:30:08 UTC 2023 refresh)
$ws.Range("D2").Value = "42.880.17"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "2.277.19"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Formula = "'249.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").Formula = "'0.641"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").Formula = "'78.76"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.23%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Formula = "'0.645"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("D10").Formula = "'41.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.60%  "
$ws.Range("D11").Formula = "'0.0973"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Formula = "'7.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "2.619.30"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Formula = "'15.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Formula = "'0.871"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("D17").Value = "2.279.77"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "42.765.72"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "0.0₃0996"
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("D21").Formula = "'72.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").Formula = "'233.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").Formula = "'2.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").Formula = "'3.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Formula = "'11.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.83%  "
$ws.Range("E27").Value = "  -4.64%  "
$ws.Range("E28").Value = "  +2.21%  "
$ws.Range("D29").Formula = "'168.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").Formula = "'20.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("D31").Formula = "'6.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.28%  "
$ws.Range("D32").Formula = "'0.0856"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.92%  "
$ws.Range("E33").Value = "  -4.22%  "
$ws.Range("D34").Formula = "'30.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.07%  "
$ws.Range("D35").Formula = "'0.128"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").Formula = "'4.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.44%  "
$ws.Range("D37").Formula = "'4.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").Formula = "'0.0305"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").Formula = "'13.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.88%  "
$ws.Range("E40").Value = "  -3.10%  "
$ws.Range("D41").Formula = "'5.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("D42").Formula = "'114.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +18.05%  "
$ws.Range("D43").Formula = "'0.209"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.53%  "
$ws.Range("D44").Formula = "'61.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Formula = "'8.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.11%  "
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").Formula = "'4.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.41%  "
$ws.Range("B48").Value = "BinanceUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D48").Formula = "'1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("E49").Value = "  -3.03%  "
$ws.Range("D50").Formula = "'1.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.49%  "
$ws.Range("D51").Formula = "'4.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.02%  "
